# separate dept from affiliations
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "PI hours": split the old combined "dept" column (which actually
#    held the full affiliation list, e.g. "['ME', 'AE', 'CSL']") into a
#    single primary "dept" column plus a new "app" (affiliations) column
#    that keeps the original list values.
# ---------------------------------------------------------------------
$piSheet = $wb.Worksheets.Item("PI hours")

# New column F = "app", formatted like the existing header cells.
$piSheet.Range("E1").Copy()
$piSheet.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$piSheet.Range("F1").Value = "app"

# Move the old full-affiliation values into the new F column...
$piSheet.Range("F2").Value = "['ME', 'AE', 'CSL']"
$piSheet.Range("F3").Value = "['ECE', 'CSL']"
$piSheet.Range("F4").Value = "['CS', 'CSL']"
$piSheet.Range("F5").Value = "['ME', 'CSL']"
$piSheet.Range("F6").Value = "['ECE', 'CSL']"

# ...and replace E with just the primary department.
$piSheet.Range("E2").Value = "ME"
$piSheet.Range("E3").Value = "ECE"
$piSheet.Range("E4").Value = "CS"
$piSheet.Range("E5").Value = "ME"
$piSheet.Range("E6").Value = "ECE"

# ---------------------------------------------------------------------
# 2. Duplicate the existing "dept hours" sheet (unchanged data) and
#    rename the copy to "unit(accumulative) hours"; it becomes the new
#    3rd sheet.
# ---------------------------------------------------------------------
$deptSheet = $wb.Worksheets.Item("dept hours")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$deptSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$unitSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$unitSheet.Name = "unit(accumulative) hours"
$unitSheet.Range("B1").Value = "unit(accumulative)"

# ---------------------------------------------------------------------
# 3. Turn the original "dept hours" sheet into the new "department
#    hours" sheet: aggregate hours/percentage by the new single "dept"
#    column from "PI hours" (ME, ECE, CS) and drop the now-unused rows.
# ---------------------------------------------------------------------
$deptSheet.Name = "department hours"

$deptSheet.Range("B2").Value = "ME"
$deptSheet.Range("C2").Value = 39
$deptSheet.Range("D2").Value = 52.34899328859061

$deptSheet.Range("B3").Value = "ECE"
$deptSheet.Range("C3").Value = 21.5
$deptSheet.Range("D3").Value = 28.85906040268456

$deptSheet.Range("B4").Value = "CS"
$deptSheet.Range("C4").Value = 14
$deptSheet.Range("D4").Value = 18.79194630872483

$deptSheet.Range("A5:D6").Delete(-4162)  # xlShiftUp
